$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "multi"
$ws1.Range("J2").Value = "sub"
$ws1.Range("A3").Value = 10
$ws1.Range("J3").Value = 4

$ws1.Range("J2").Select()
